# Update the edges worksheet: rename header columns and rewrite the
# full Parent_Node/Child_Node edge list (62 rows incl. header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Parent_Node", "Child_Node"),
    @("FIT101.Pv", "MV101.Status"),
    @("MV101.Status", "P401.Status"),
    @("P101.Status", "MV201.Status"),
    @("P101.Status", "P3_STATE"),
    @("P2_STATE", "P1_STATE"),
    @("MV201.Status", "P203.Status"),
    @("P201.Status", "P205.Status"),
    @("P201.Status", "AIT201.Pv"),
    @("P203.Status", "P201.Status"),
    @("P203.Status", "FIT201.Pv"),
    @("P205.Status", "AIT303.Pv"),
    @("P3_STATE", "MV302.Status"),
    @("P3_STATE", "MV304.Status"),
    @("P3_STATE", "MV303.Status"),
    @("P3_STATE", "MV301.Status"),
    @("DPIT301.Pv", "P301.Status"),
    @("MV301.Status", "P602.Status"),
    @("P301.Status", "P3_STATE"),
    @("P301.Status", "FIT301.Pv"),
    @("P301.Status", "AIT503.Pv"),
    @("P301.Status", "AIT302.Pv"),
    @("P301.Status", "PIT503.Pv"),
    @("P301.Status", "PIT502.Pv"),
    @("P301.Status", "AIT501.Pv"),
    @("P301.Status", "PIT501.Pv"),
    @("P301.Status", "FIT601.Pv"),
    @("P301.Status", "P1_STATE"),
    @("P301.Status", "P101.Status"),
    @("P301.Status", "MV201.Status"),
    @("P401.Status", "UV401.Status"),
    @("P401.Status", "P5_STATE"),
    @("UV401.Status", "FIT503.Pv"),
    @("UV401.Status", "AIT504.Pv"),
    @("P5_STATE", "MV502.Status"),
    @("P5_STATE", "MV504.Status"),
    @("P5_STATE", "MV503.Status"),
    @("P5_STATE", "P6_STATE"),
    @("FIT503.Pv", "P501.Status"),
    @("AIT503.Pv", "P101.Status"),
    @("AIT503.Pv", "LSH601.Alarm"),
    @("AIT503.Pv", "P601.Status"),
    @("PIT502.Pv", "MV501.Status"),
    @("P501.Status", "MV502.Status"),
    @("P501.Status", "FIT401.Pv"),
    @("P501.Status", "PIT501.Pv"),
    @("P501.Status", "PIT503.Pv"),
    @("P501.Status", "FIT501.Pv"),
    @("P501.Status", "FIT502.Pv"),
    @("P501.Status", "AIT503.Pv"),
    @("P501.Status", "AIT502.Pv"),
    @("P501.Status", "DPIT301.Pv"),
    @("P501.Status", "P3_STATE"),
    @("P501.Status", "AIT301.Pv"),
    @("MV501.Status", "P5_STATE"),
    @("P6_STATE", "P2_STATE"),
    @("P601.Status", "P1_STATE"),
    @("LSH601.Alarm", "PIT502.Pv"),
    @("LSH601.Alarm", "P3_STATE"),
    @("LSH601.Alarm", "AIT504.Pv"),
    @("LSH601.Alarm", "P1_STATE"),
    @("LSH601.Alarm", "P203.Status")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
